$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Civ")
$ws.Range("B1:B17").Select()
Write-Host "done"
